$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '65.774.46'
$ws.Range('E2').Value = '  -4.44%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.271.64'
$ws.Range('E3').Value = '  -4.83%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '554.61'
$ws.Range('E5').Value = '  -3.28%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '183.59'
$ws.Range('E6').Value = '  -2.78%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.05%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.592'

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '3.262.57'

# Row 10
$ws.Range('E10').Value = '  -7.37%  '

# Row 11
$ws.Range('E11').Value = '  -4.26%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '47.25'
$ws.Range('E12').Value = '  -7.13%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000266'
$ws.Range('E13').Value = '  -5.55%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '8.57'
$ws.Range('E14').Value = '  -4.99%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '625.94'
$ws.Range('E15').Value = '  -1.40%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.795.20'
$ws.Range('E16').Value = '  -4.68%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '65.610.78'
$ws.Range('E17').Value = '  -4.42%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '17.85'
$ws.Range('E18').Value = '  -0.60%  '

# Row 19
$ws.Range('E19').Value = '  -3.20%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '3.277.28'
$ws.Range('E20').Value = '  -4.52%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.35'
$ws.Range('E21').Value = '  -6.36%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.903'
$ws.Range('E22').Value = '  -3.27%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '17.61'
$ws.Range('E23').Value = '  -0.06%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '106.10'
$ws.Range('E24').Value = '  +8.52%  '

# Row 25
$ws.Range('E25').Value = '  -6.77%  '

# Row 26
$ws.Range('E26').Value = '  -6.34%  '

# Row 27
$ws.Range('E27').Value = '  -5.71%  '

# Row 28
$ws.Range('E28').Value = '  -2.58%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '8.65'
$ws.Range('E29').Value = '  -5.30%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '30.38'
$ws.Range('E30').Value = '  -5.17%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.02'
$ws.Range('E31').Value = '  -5.60%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.25'

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '10.99'
$ws.Range('E33').Value = '  -4.09%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '540.63'
$ws.Range('E34').Value = '  +10.65%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.105'
$ws.Range('E35').Value = '  -3.08%  '

# Row 36
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '57.27'
$ws.Range('E36').Value = '  -5.75%  '

# Row 37
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  -0.12%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.667.30'
$ws.Range('E38').Value = '  +0.69%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.39'
$ws.Range('E39').Value = '  -0.70%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.0₃0727'
$ws.Range('E40').Value = '  -6.57%  '

# Row 41
$ws.Range('E41').Value = '  -1.22%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.72'
$ws.Range('E42').Value = '  -5.00%  '

# Row 43
$ws.Range('B43').Value = 'CoreDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.36'
$ws.Range('E43').Value = '  -4.85%  '

# Row 44
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '32.49'
$ws.Range('E44').Value = '  -4.04%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.336'
$ws.Range('E45').Value = '  -7.87%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.25'
$ws.Range('E46').Value = '  -2.17%  '

# Row 47
$ws.Range('E47').Value = '  -4.69%  '

# Row 48
$ws.Range('E48').Value = '  -6.03%  '

# Row 49
$ws.Range('E49').Value = '  -3.21%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.997'
$ws.Range('E50').Value = '  -0.12%  '

# Row 51
$ws.Range('E51').Value = '  +2.20%  '

